# "Generate Report for Handback"
#
# This script regenerates the localization-status report after a handback:
#  - Updates the Status of both locales (zh-cn, de-de) from "Ready for handoff"
#    to "Handed back: in sync with en-US" on the Overview sheet and on each
#    locale's detail sheet.
#  - Refreshes the "Latest Handback DateTime" timestamps for both locales.
#  - Clears the stale "Error Detail" (out-of-date handback warning) now that
#    the handback is in sync.
#  - Re-sizes a few columns to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for zh-cn (E) and de-de (F) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-24 10:21:00"
$zhcn.Range("K3").Value = "2016-10-24 10:21:00"
$zhcn.Range("P3").Value = ""

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-10-24 10:21:17"
$dede.Range("K3").Value = "2016-10-24 10:21:17"
$dede.Range("P3").Value = ""

# --- Column width adjustments to fit the new content ---
$overview.Range("E1").ColumnWidth = 29.15
$overview.Range("F1").ColumnWidth = 29.15

$zhcn.Range("C1").ColumnWidth = 29.15
$zhcn.Range("P1").ColumnWidth = 12.8

$dede.Range("C1").ColumnWidth = 29.15
$dede.Range("P1").ColumnWidth = 12.8
